$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Cells.Item(7, 1).Value = "refrigerator_text"
$ws.Cells.Item(7, 2).Value = "text"

$ws.Rows.Item(7).RowHeight = 12.75

$ws.Range("B7").Select()
